$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I, shifting Celular/CPF/CEP right by one
$ws.Columns("I:I").Insert()

# Ensure phone/doc columns stay formatted as text (avoid numeric auto-conversion)
$ws.Range("I1:L11").NumberFormat = "@"

# Set header for new column and adjust header for shifted Celular column
$ws.Range("I1").Value = "Telefone"
$ws.Range("J1").Value = "Celular"

# Row 2
$ws.Range("B2").Value = "Gustavo Ferreira Santana"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = "Masculino"
$ws.Range("E2").Value = "Branco"
$ws.Range("F2").Value = "Sueli Borges Santana Ferreira"
$ws.Range("G2").Value = "Em idade escolar."
$ws.Range("H2").Value = "Fora da força de trabalho"
$ws.Range("I2").Value = "53 3566-8652"
$ws.Range("J2").Value = "53 91088-6697"
$ws.Range("K2").Value = "934.116.140-13"
$ws.Range("L2").Value = "292552690"

# Row 3
$ws.Range("B3").Value = "Roberto Gustavo Santos Ribeira"
$ws.Range("C3").Value = 19
$ws.Range("D3").Value = "Masculino"
$ws.Range("E3").Value = "Pardo"
$ws.Range("F3").Value = "Carolina Simone Santos Gustavo"
$ws.Range("G3").Value = "Em idade escolar."
$ws.Range("H3").Value = "Empregado: Trabalhador doméstico (sem CLT)"
$ws.Range("I3").Value = "97 0757-5038"
$ws.Range("J3").Value = "97 91758-7448"
$ws.Range("K3").Value = "952.836.470-56"
$ws.Range("L3").Value = "645356793"

# Row 4
$ws.Range("B4").Value = "Larissa Tatiana Martins"
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = "Feminino"
$ws.Range("E4").Value = "Outro"
$ws.Range("F4").Value = "Flávia Andressa Martins Tatiana"
$ws.Range("G4").Value = "Em idade escolar."
$ws.Range("H4").Value = "Fora da força de trabalho"
$ws.Range("I4").Value = "88 2239-1378"
$ws.Range("J4").Value = "88 93770-1262"
$ws.Range("K4").Value = "592.811.940-23"
$ws.Range("L4").Value = "660156060"

# Row 5
$ws.Range("B5").Value = "Sérgio Cavalcanti Barros Ribeira Pires"
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = "Masculino"
$ws.Range("E5").Value = "Branco"
$ws.Range("F5").Value = "Rafaela Barros Cavalcanti"
$ws.Range("G5").Value = "Sem instrução"
$ws.Range("H5").Value = "Empregado: Setor privado (CLT)"
$ws.Range("I5").Value = "79 4723-9692"
$ws.Range("J5").Value = "79 93616-1771"
$ws.Range("K5").Value = "591.100.060-10"
$ws.Range("L5").Value = "430784099"

# Row 6
$ws.Range("B6").Value = "Fernando Marcelo Martins"
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = "Masculino"
$ws.Range("E6").Value = "Pardo"
$ws.Range("F6").Value = "Eduarda Pires Martins Marcelo"
$ws.Range("G6").Value = "Em idade escolar."
$ws.Range("H6").Value = "Fora da força de trabalho"
$ws.Range("I6").Value = "34 8557-2628"
$ws.Range("J6").Value = "34 98097-9150"
$ws.Range("K6").Value = "542.851.200-89"
$ws.Range("L6").Value = "096772649"

# Row 7
$ws.Range("B7").Value = "Samuel Carvalho"
$ws.Range("C7").Value = 24
$ws.Range("D7").Value = "Masculino"
$ws.Range("E7").Value = "Branco"
$ws.Range("F7").Value = "Lorena Pereira Carvalho"
$ws.Range("G7").Value = "Em idade escolar."
$ws.Range("H7").Value = "Empregado: Setor privado (sem CLT)"
$ws.Range("I7").Value = "43 6946-3963"
$ws.Range("J7").Value = "43 99291-6759"
$ws.Range("K7").Value = "113.811.910-53"
$ws.Range("L7").Value = "219400081"

# Row 8
$ws.Range("B8").Value = "Miguel Azevedo"
$ws.Range("C8").Value = 28
$ws.Range("D8").Value = "Masculino"
$ws.Range("E8").Value = "Pardo"
$ws.Range("F8").Value = "Isabela Karla Azevedo"
$ws.Range("G8").Value = "Médio completo"
$ws.Range("H8").Value = "Empregado: Setor publico (sem CLT)"
$ws.Range("I8").Value = "12 5288-6779"
$ws.Range("J8").Value = "12 96640-9757"
$ws.Range("K8").Value = "241.821.300-57"
$ws.Range("L8").Value = "064478034"

# Row 9
$ws.Range("B9").Value = "Márcia Santana Silveira Ribeira"
$ws.Range("C9").Value = 18
$ws.Range("D9").Value = "Feminino"
$ws.Range("E9").Value = "Parda"
$ws.Range("F9").Value = "Sueli Ribeira Silveira Santana"
$ws.Range("G9").Value = "Em idade escolar."
$ws.Range("H9").Value = "Empregado: Setor privado (CLT)"
$ws.Range("I9").Value = "55 3579-4968"
$ws.Range("J9").Value = "55 90199-5757"
$ws.Range("K9").Value = "572.447.780-74"
$ws.Range("L9").Value = "049526582"

# Row 10
$ws.Range("B10").Value = "Alexandre Mendes"
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "Masculino"
$ws.Range("E10").Value = "Branco"
$ws.Range("F10").Value = "Jessica Mendes"
$ws.Range("G10").Value = "Em idade escolar."
$ws.Range("H10").Value = "Fora da força de trabalho"
$ws.Range("I10").Value = "54 0394-6492"
$ws.Range("J10").Value = "54 92238-3698"
$ws.Range("K10").Value = "034.705.200-20"
$ws.Range("L10").Value = "808547634"

# Row 11
$ws.Range("B11").Value = "Vinícius Gonçalves Rodrigues Oliveira"
$ws.Range("C11").Value = 13
$ws.Range("D11").Value = "Masculino"
$ws.Range("E11").Value = "Branco"
$ws.Range("F11").Value = "Tatiana Gonçalves Oliveira Rodrigues"
$ws.Range("G11").Value = "Em idade escolar."
$ws.Range("H11").Value = "Fora da força de trabalho"
$ws.Range("I11").Value = "69 8126-8908"
$ws.Range("J11").Value = "69 96346-3282"
$ws.Range("K11").Value = "459.864.800-74"
$ws.Range("L11").Value = "224266837"
